# Update test case content: the site/app dropped "Collection" as a
# separate nav item, so the test case wording was revised to only
# reference Home / Browse / Add.
#
# Cell values are set in the same order the author re-typed them so the
# shared-string table is (re)built in the matching order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "1) Click on tab link to`n(Home, Browse, Add)"
$ws.Range("B5").Value = "User at Home Content Page "
$ws.Range("B7").Value = "User at Bottom Home Content Page (Browse,  Add)"
$ws.Range("C6").Value = "1) Clicking will link user to browse page"
$ws.Range("B6").Value = "Browse tab link (centre browse button)"
$ws.Range("B8").Value = " Browse or Add button link"

# The author's last selection before saving was cell B8.
$ws.Range("B8").Select() | Out-Null
